$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handback" - populate the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) now that
# the localized files have come back in sync with en-US, and flip the
# Status column from "Ready for handoff" to "Handed back: in sync with
# en-US" for both locale report sheets.
# ---------------------------------------------------------------------------

$zhRowA = @{
    2 = @{ Base = "310c4156-6a65-421f-a2b7-76c22e5e8b5b"; Md = "310c4156-6a65-421f-a2b7-76c22e5e8b5b.md"; Xlf = "310c4156-6a65-421f-a2b7-76c22e5e8b5b.7fc0c9500424b0b8aa7497387d8b4cd441b5a7f6.zh-cn.xlf" }
    3 = @{ Base = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5"; Md = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.md"; Xlf = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.9f01dbcb58abff7f51af3e5c0bcb7c0674bc54e1.zh-cn.xlf" }
}

$zhMdUrl = @{
    2 = "https://github.com/OpenLocalizationTest/oltest/blob/4e09b2f77302c720b34895735ef62e9db29af3b6/e2e/310c4156-6a65-421f-a2b7-76c22e5e8b5b.md"
    3 = "https://github.com/OpenLocalizationTest/oltest/blob/4e09b2f77302c720b34895735ef62e9db29af3b6/e2e/c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.md"
}
$zhXlfUrl = @{
    2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd5be9d684093e80dca9707425ae9116522331f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/310c4156-6a65-421f-a2b7-76c22e5e8b5b.7fc0c9500424b0b8aa7497387d8b4cd441b5a7f6.zh-cn.xlf"
    3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd5be9d684093e80dca9707425ae9116522331f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.9f01dbcb58abff7f51af3e5c0bcb7c0674bc54e1.zh-cn.xlf"
}

$deRowA = @{
    2 = @{ Base = "310c4156-6a65-421f-a2b7-76c22e5e8b5b"; Md = "310c4156-6a65-421f-a2b7-76c22e5e8b5b.md"; Xlf = "310c4156-6a65-421f-a2b7-76c22e5e8b5b.7fc0c9500424b0b8aa7497387d8b4cd441b5a7f6.de-de.xlf" }
    3 = @{ Base = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5"; Md = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.md"; Xlf = "c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.9f01dbcb58abff7f51af3e5c0bcb7c0674bc54e1.de-de.xlf" }
}

$deMdUrl = @{
    2 = "https://github.com/OpenLocalizationTest/oltest/blob/4e09b2f77302c720b34895735ef62e9db29af3b6/e2e/310c4156-6a65-421f-a2b7-76c22e5e8b5b.md"
    3 = "https://github.com/OpenLocalizationTest/oltest/blob/4e09b2f77302c720b34895735ef62e9db29af3b6/e2e/c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.md"
}
$deXlfUrl = @{
    2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb339cffd4d653840faca28a9b8741909969115e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/310c4156-6a65-421f-a2b7-76c22e5e8b5b.7fc0c9500424b0b8aa7497387d8b4cd441b5a7f6.de-de.xlf"
    3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb339cffd4d653840faca28a9b8741909969115e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c237d730-e9e1-4d92-aa9b-6e5ba18c56e5.9f01dbcb58abff7f51af3e5c0bcb7c0674bc54e1.de-de.xlf"
}

# ---- Overview sheet --------------------------------------------------------
# The Overview sheet mirrors the same "Status" shared string for each
# locale column, so it needs to move in lock-step with the per-locale
# report sheets below.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $wsOverview.Range("B$row").Value = "Handed back: in sync with en-US"
    $wsOverview.Range("C$row").Value = "Handed back: in sync with en-US"
}

# ---- zh-cn sheet ----------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

foreach ($row in 2, 3) {
    $ws.Range("C$row").Value = "Handed back: in sync with en-US"
    $ws.Range("H$row").Value = "2016-03-24 13:22:45"
}

foreach ($row in 2, 3) {
    $info = $zhRowA[$row]
    $ws.Hyperlinks.Add($ws.Range("F$row"), $zhMdUrl[$row], "", "", $info.Md)
    $ws.Hyperlinks.Add($ws.Range("G$row"), $zhXlfUrl[$row], "", "", $info.Xlf)
}

# ---- de-de sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

foreach ($row in 2, 3) {
    $ws.Range("C$row").Value = "Handed back: in sync with en-US"
    $ws.Range("H$row").Value = "2016-03-24 13:22:52"
}

foreach ($row in 2, 3) {
    $info = $deRowA[$row]
    $ws.Hyperlinks.Add($ws.Range("F$row"), $deMdUrl[$row], "", "", $info.Md)
    $ws.Hyperlinks.Add($ws.Range("G$row"), $deXlfUrl[$row], "", "", $info.Xlf)
}
